# Auto-generated edit script: updates cryptos.xlsx price/volume data
# per commit "Updated cryptos list on Sun Feb 18 05:45:50 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '51.644.16'
$ws.Cells.Item(2, 5).Value = '  -0.46%  '
$ws.Cells.Item(3, 4).Value = '2.802.41'
$ws.Cells.Item(3, 5).Value = '  +0.66%  '
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$ws.Cells.Item(5, 4).Value = "'" + '355.96'
$ws.Cells.Item(5, 5).Value = '  -0.27%  '
$ws.Cells.Item(6, 4).Value = "'" + '109.58'
$ws.Cells.Item(6, 5).Value = '  +0.25%  '
$ws.Cells.Item(7, 4).Value = "'" + '0.556'
$ws.Cells.Item(7, 5).Value = '  -0.29%  '
$ws.Cells.Item(8, 4).Value = "'" + '1.00'
$ws.Cells.Item(8, 5).Value = '  +0.03%  '
$ws.Cells.Item(9, 4).Value = "'" + '0.623'
$ws.Cells.Item(9, 5).Value = '  +5.53%  '
$ws.Cells.Item(10, 4).Value = "'" + '40.23'
$ws.Cells.Item(10, 5).Value = '  -0.17%  '
$ws.Cells.Item(11, 5).Value = '  +1.23%  '
$ws.Cells.Item(12, 5).Value = '  -1.05%  '
$ws.Cells.Item(13, 4).Value = "'" + '20.04'
$ws.Cells.Item(13, 5).Value = '  +2.94%  '
$ws.Cells.Item(14, 4).Value = "'" + '7.80'
$ws.Cells.Item(14, 5).Value = '  +3.10%  '
$ws.Cells.Item(15, 4).Value = '3.239.45'
$ws.Cells.Item(15, 5).Value = '  +0.46%  '
$ws.Cells.Item(16, 4).Value = '2.797.99'
$ws.Cells.Item(16, 5).Value = '  +0.82%  '
$ws.Cells.Item(17, 4).Value = "'" + '0.943'
$ws.Cells.Item(17, 5).Value = '  -0.05%  '
$ws.Cells.Item(18, 4).Value = '51.651.00'
$ws.Cells.Item(18, 5).Value = '  -0.33%  '
$ws.Cells.Item(19, 4).Value = "'" + '7.77'
$ws.Cells.Item(20, 4).Value = "'" + '3.19'
$ws.Cells.Item(20, 5).Value = '  +3.10%  '
$ws.Cells.Item(21, 4).Value = "'" + '13.41'
$ws.Cells.Item(21, 5).Value = '  +2.22%  '
$ws.Cells.Item(22, 4).Value = '0.0₃0973'
$ws.Cells.Item(22, 5).Value = '  -0.35%  '
$ws.Cells.Item(23, 5).Value = '  +0.62%  '
$ws.Cells.Item(24, 4).Value = "'" + '268.36'
$ws.Cells.Item(24, 5).Value = '  -0.52%  '
$ws.Cells.Item(25, 5).Value = '  +1.07%  '
$ws.Cells.Item(26, 4).Value = "'" + '0.999'
$ws.Cells.Item(26, 5).Value = '  +0.01%  '
$ws.Cells.Item(27, 4).Value = "'" + '26.12'
$ws.Cells.Item(27, 5).Value = '  -1.39%  '
$ws.Cells.Item(28, 5).Value = '  +1.07%  '
$ws.Cells.Item(29, 4).Value = "'" + '10.37'
$ws.Cells.Item(29, 5).Value = '  +0.59%  '
$ws.Cells.Item(30, 4).Value = "'" + '37.49'
$ws.Cells.Item(30, 5).Value = '  +8.94%  '
$ws.Cells.Item(31, 5).Value = '  +4.71%  '
$ws.Cells.Item(32, 4).Value = "'" + '6.40'
$ws.Cells.Item(32, 5).Value = '  +11.71%  '
$ws.Cells.Item(33, 4).Value = "'" + '52.17'
$ws.Cells.Item(33, 5).Value = '  +0.07%  '
$ws.Cells.Item(34, 4).Value = "'" + '5.67'
$ws.Cells.Item(34, 5).Value = '  +9.11%  '
$ws.Cells.Item(35, 4).Value = "'" + '0.0447'
$ws.Cells.Item(35, 5).Value = '  -4.80%  '
$ws.Cells.Item(36, 4).Value = "'" + '0.0857'
$ws.Cells.Item(36, 5).Value = '  +1.42%  '
$ws.Cells.Item(37, 4).Value = "'" + '1.00'
$ws.Cells.Item(37, 5).Value = '  -0.01%  '
$ws.Cells.Item(38, 4).Value = "'" + '18.93'
$ws.Cells.Item(38, 5).Value = '  +0.15%  '
$ws.Cells.Item(39, 5).Value = '  -1.59%  '
$ws.Cells.Item(40, 5).Value = '  +0.25%  '
$ws.Cells.Item(41, 5).Value = '  +0.33%  '
$ws.Cells.Item(42, 5).Value = '  -4.97%  '
$ws.Cells.Item(43, 4).Value = "'" + '119.94'
$ws.Cells.Item(43, 5).Value = '  +0.08%  '
$ws.Cells.Item(44, 2).Value = 'EnergySwap'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(44, 4).Value = "'" + '21.91'
$ws.Cells.Item(44, 5).Value = '  +0.36%  '
$ws.Cells.Item(45, 2).Value = 'WEMIXToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(45, 4).Value = "'" + '2.19'
$ws.Cells.Item(45, 5).Value = '  -2.41%  '
$ws.Cells.Item(46, 2).Value = 'Maker'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(46, 4).Value = '2.141.34'
$ws.Cells.Item(46, 5).Value = '  +2.59%  '
$ws.Cells.Item(47, 2).Value = 'NEARProtocol'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(47, 4).Value = "'" + '3.44'
$ws.Cells.Item(47, 5).Value = '  +5.39%  '
$ws.Cells.Item(48, 5).Value = '  +7.24%  '
$ws.Cells.Item(49, 4).Value = "'" + '0.926'
$ws.Cells.Item(49, 5).Value = '  -3.50%  '
$ws.Cells.Item(50, 5).Value = '  +10.74%  '
$ws.Cells.Item(51, 4).Value = "'" + '0.221'
$ws.Cells.Item(51, 5).Value = '  +16.72%  '
